$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Normalize "de"/"el"/"los" -> "De"/"El"/"Los" in specific cells
$ws.Range("A8").Value = "Ciudad De México"
$ws.Range("A14").Value = "Estado De México"
$ws.Range("B14").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B24").Value = "San Juan De Los Lagos"
$ws.Range("B25").Value = "San Miguel El Alto"
$ws.Range("B26").Value = "Tepatitlán De Morelos"
$ws.Range("B29").Value = "Valle De Guadalupe"

# Delete the trailing metadata rows (51-55), shifting cells up
$ws.Range("A51:A55").EntireRow.Delete()
